$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.032.90"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "1.666.99"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5092"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2648"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06383"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07434"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").Value = "1.666.88"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.507"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5814"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.60%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.11%  "

$ws.Range("D17").Value = "26.075.54"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.928"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.20%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.180"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.603"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1198"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06589"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.47%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.337"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("E30").Value = "  -1.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.543"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.507"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.654"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.016"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6124"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.88%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.685"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.313"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.93%  "

$ws.Range("D39").Value = "1.094.53"
$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01598"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8702"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "

$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").Value = "1.816.81"
$ws.Range("E44").Value = "  -1.75%  "

$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.007"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.044"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05229"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4287"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.026"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.46%  "
